$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H42").Value = 1470
$ws.Range("I42").Value = 1724.2
$ws.Range("J42").Value = 199
$ws.Range("K42").Value = 5172.6
$ws.Range("L42").Value = 597
$ws.Range("M42").Value = -4942.6
$ws.Range("N42").Value = -1057

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 2930.4285
$ws.Range("I70").Value = 1002
$ws.Range("J70").Value = 3078.7693
$ws.Range("K70").Value = 3006
$ws.Range("L70").Value = 9236.3079
$ws.Range("M70").Value = -2736
$ws.Range("N70").Value = -9776.3079

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H73").Value = 2930.4285
$ws.Range("I73").Value = 1002
$ws.Range("J73").Value = 3078.7693
$ws.Range("K73").Value = 3006
$ws.Range("L73").Value = 9236.3079
$ws.Range("M73").Value = -2070
$ws.Range("N73").Value = -11108.3079

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 2600.2273
$ws.Range("I100").Value = 2541.1765
$ws.Range("J100").Value = 2801
$ws.Range("K100").Value = 2541.1765
$ws.Range("L100").Value = 2801
$ws.Range("M100").Value = -2000.1765
$ws.Range("N100").Value = -3883

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H127").Value = 2357.4546
$ws.Range("I127").Value = 2399.7896
$ws.Range("J127").Value = 2089.3333
$ws.Range("K127").Value = 7199.3688
$ws.Range("L127").Value = 6267.999899999999
$ws.Range("M127").Value = -2239.3688
$ws.Range("N127").Value = -16187.9999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H131").Value = 7329.2974
$ws.Range("I131").Value = 4928.577
$ws.Range("K131").Value = 14785.731
$ws.Range("M131").Value = -9745.731

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 4542.25
$ws.Range("I132").Value = 4542.25
$ws.Range("K132").Value = 13626.75
$ws.Range("M132").Value = -11096.75

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 121171.57
$ws.Range("J137").Value = 240953.08
$ws.Range("L137").Value = 722859.24
$ws.Range("N137").Value = -727959.24

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2630.7693
$ws.Range("I138").Value = 1781.9166
$ws.Range("J138").Value = 3008.037
$ws.Range("K138").Value = 5345.7498
$ws.Range("L138").Value = 9024.110999999999
$ws.Range("M138").Value = -205.7497999999996
$ws.Range("N138").Value = -19304.111

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3342.7437
$ws.Range("I61").Value = 2569.6875
$ws.Range("K61").Value = 2569.6875
$ws.Range("M61").Value = -2357.6875

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 38148.594
$ws.Range("I74").Value = 26869.742
$ws.Range("J74").Value = 66345.71000000001
$ws.Range("K74").Value = 26869.742
$ws.Range("L74").Value = 66345.71000000001
$ws.Range("M74").Value = -25995.742
$ws.Range("N74").Value = -68093.71000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 38148.594
$ws.Range("I77").Value = 26869.742
$ws.Range("J77").Value = 66345.71000000001
$ws.Range("K77").Value = 134348.71
$ws.Range("L77").Value = 331728.55
$ws.Range("M77").Value = -129980.71
$ws.Range("N77").Value = -340464.55

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 3342.7437
$ws.Range("I136").Value = 2569.6875
$ws.Range("K136").Value = 7709.0625
$ws.Range("M136").Value = -5159.0625

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2046.2
$ws.Range("I99").Value = 2115.4
$ws.Range("J99").Value = 1838.6
$ws.Range("K99").Value = 2115.4
$ws.Range("L99").Value = 1838.6
$ws.Range("M99").Value = -617.4000000000001
$ws.Range("N99").Value = -4834.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1926.0209
$ws.Range("I134").Value = 1782.9762
$ws.Range("J134").Value = 2927.3333
$ws.Range("K134").Value = 5348.9286
$ws.Range("L134").Value = 8781.999899999999
$ws.Range("M134").Value = -2813.9286
$ws.Range("N134").Value = -13851.9999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 16855.428

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H89").Value = 16855.428

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 4892.8213
$ws.Range("I99").Value = 4558.55
$ws.Range("J99").Value = 5728.5
$ws.Range("K99").Value = 4558.55
$ws.Range("L99").Value = 5728.5
$ws.Range("M99").Value = -3060.55
$ws.Range("N99").Value = -8724.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 4892.8213
$ws.Range("I126").Value = 4558.55
$ws.Range("J126").Value = 5728.5
$ws.Range("K126").Value = 13675.65
$ws.Range("L126").Value = 17185.5
$ws.Range("M126").Value = -11205.65
$ws.Range("N126").Value = -22125.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 5152.7856
$ws.Range("I132").Value = 4512.636
$ws.Range("K132").Value = 13537.908
$ws.Range("M132").Value = -11007.908

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 4187.7383
$ws.Range("I134").Value = 4235.0625
$ws.Range("K134").Value = 12705.1875
$ws.Range("M134").Value = -10170.1875

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 28478.715
$ws.Range("J80").Value = 25130.3
$ws.Range("L80").Value = 75390.89999999999
$ws.Range("N80").Value = -77262.89999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H83").Value = 28478.715
$ws.Range("J83").Value = 25130.3
$ws.Range("L83").Value = 226172.7
$ws.Range("N83").Value = -235532.7

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 1540.4445
$ws.Range("J107").Value = 1884.2858
$ws.Range("L107").Value = 5652.857400000001
$ws.Range("N107").Value = -9492.857400000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 1850.5555
$ws.Range("I137").Value = 1185.1
$ws.Range("J137").Value = 2682.375
$ws.Range("K137").Value = 3555.3
$ws.Range("L137").Value = 8047.125
$ws.Range("M137").Value = 1544.7
$ws.Range("N137").Value = -18247.125

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 24679.9
$ws.Range("I136").Value = 1758.6451
$ws.Range("J136").Value = 103630.89
$ws.Range("K136").Value = 5275.9353
$ws.Range("L136").Value = 310892.67
$ws.Range("M136").Value = -2725.9353

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 10288.111
$ws.Range("I62").Value = 7219
$ws.Range("K62").Value = 7219
$ws.Range("M62").Value = -6595

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H65").Value = 10288.111
$ws.Range("I65").Value = 7219
$ws.Range("K65").Value = 36095
$ws.Range("M65").Value = -32975

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 407.45834
$ws.Range("I100").Value = 331.46667
$ws.Range("J100").Value = 534.1111
$ws.Range("K100").Value = 662.93334
$ws.Range("L100").Value = 1068.2222
$ws.Range("M100").Value = -121.93334
$ws.Range("N100").Value = -2150.2222

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H125").Value = 63996.4
$ws.Range("J125").Value = 63996.4
$ws.Range("L125").Value = 63996.4
$ws.Range("N125").Value = -73836.39999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 44331.234
$ws.Range("J136").Value = 56240.35
$ws.Range("L136").Value = 168721.05
$ws.Range("N136").Value = -173821.05
